$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update date / shift ---
$ws.Range("A2").Value = "2025-02-22T18:36"
$ws.Range("B2").Value = "SHIFT_2"

# --- Row 2: the "*_cnt" numeric columns become the literal text "0" ---
# (leading apostrophe forces Excel to store a text value instead of
#  re-interpreting the digit string as a number)
$cntCols = @("D","F","H","J","L","N","P","R","T","V","X","Z","AB","AD","AF")
foreach ($col in $cntCols) {
    $ws.Range($col + "2").Value = "'0"
}

# --- Row 2: Remark / QA-Sign / Engg-Sign become free text ---
$ws.Range("AG2").Value = "testing"
$ws.Range("AH2").Value = "suriya"
$ws.Range("AI2").Value = "suriya"

# --- Remove the old row 3 entirely (dimension shrinks to A1:AI2) ---
$ws.Rows(3).Delete()
